# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# "OFF" sheet (sheet1) - row 2 (label "H")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 218
$wsOff.Range("C2").Value = 144
$wsOff.Range("D2").Value = 24
$wsOff.Range("E2").Value = 8
$wsOff.Range("F2").Value = 3
$wsOff.Range("G2").Value = 3

# "DEF" sheet (sheet2) - row 2 (label "H")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 208
$wsDef.Range("C2").Value = 149
$wsDef.Range("D2").Value = 48
